$d = $word.ActiveDocument

# --- Step 1: fix the title text --------------------------------------
# "PROJECT X: AUTOMATED ATTENDANC " (with a trailing space, truncated
# word) becomes "PROJECT X: AUTOMATED ATTENDANCE" (fixed spelling, no
# longer needs xml:space="preserve" since there's no leading/trailing
# whitespace left).
[void]$d.Content.Find.Execute("PROJECT X: AUTOMATED ATTENDANC ", $false, $false, $false, $false, $false, $true, 1, $false, "PROJECT X: AUTOMATED ATTENDANCE", 2)

# --- Step 2: relocate the "_GoBack" bookmark --------------------------
# It currently sits at the start of the second paragraph (right before
# "LOW LEVEL DIAGRAM"). It needs to move to the end of the first
# paragraph, right after the title text and before a new trailing
# space run.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# The first paragraph's Range now ends at 32 (31 chars of text + the
# paragraph mark), so position 31 is immediately after "ATTENDANCE".
# Insert a one-character placeholder there first so that position 31
# is no longer the very end of the paragraph - this avoids an engine
# quirk where adding a zero-length bookmark exactly at the paragraph
# end places its markers on the wrong side of the paragraph boundary.
$titleEnd = $d.Paragraphs(1).Range.End - 1
$placeholder = $d.Range($titleEnd, $titleEnd)
$placeholder.InsertAfter("X")

# Re-create the "_GoBack" bookmark as a collapsed range between the
# title text and the placeholder character.
$bmRange = $d.Range($titleEnd, $titleEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Turn the placeholder into the new run's actual content: a single
# space, matching the formatting of the title run.
$newRunRange = $d.Range($titleEnd, $titleEnd + 1)
$newRunRange.Text = " "
$newRunRange.Font.Bold = $true
$newRunRange.Font.Size = 20
